# Minor fix for temp files in optimisation mode.
#
# For each of Sheet1..Sheet9, replace the single "Hello world" cell with
# three cells: A1="Hello 1", A2="Hello 2", A4="Hello 3" (row 3 left blank).
# Then remove Sheet10 entirely.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le 9; $i++) {
    $ws = $wb.Worksheets.Item("Sheet$i")
    $ws.Range("A1").Value = "Hello 1"
    $ws.Range("A2").Value = "Hello 2"
    $ws.Range("A4").Value = "Hello 3"
}

$sheet10 = $wb.Worksheets.Item("Sheet10")
$sheet10.Delete()

# Deleting Sheet10 (the last/active sheet) shifts Excel's active-sheet
# selection to Sheet9; restore Sheet1 as the selected/active tab to match
# the original workbook's tabSelected state.
$wb.Worksheets.Item("Sheet1").Activate()
